$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 5 - OUAGADOUGOU station (UVM00065503)
# ---------------------------------------------------------------------------
$ws.Range("E5").Value = "OUAGADOUGOU"
$ws.Range("F5").Value = 12.35
$ws.Range("G5").Value = -1.51
$ws.Range("H5").Value = 316
$ws.Range("I5").Value = "UVM00065503"
$ws.Range("J5").Value = 65503
$ws.Range("K5").Value = 1973
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1
$ws.Range("N5").Value = 2021
$ws.Range("O5").Value = 10
$ws.Range("P5").Value = 1

# ---------------------------------------------------------------------------
# Row 6 - DONGTAI station (CHM00058251)
# ---------------------------------------------------------------------------
$ws.Range("E6").Value = "DONGTAI"
$ws.Range("F6").Value = 32.85
$ws.Range("G6").Value = 120.28
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = "CHM00058251"
$ws.Range("J6").Value = 58251
$ws.Range("K6").Value = 1953
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1
$ws.Range("N6").Value = 2021
$ws.Range("O6").Value = 10
$ws.Range("P6").Value = 1

# ---------------------------------------------------------------------------
# Row 7 - DONGTAI station (CHM00058252)
# ---------------------------------------------------------------------------
$ws.Range("E7").Value = "DONGTAI"
$ws.Range("F7").Value = 32.85
$ws.Range("G7").Value = 120.28
$ws.Range("H7").Value = 5
$ws.Range("I7").Value = "CHM00058252"
$ws.Range("J7").Value = 58251
$ws.Range("K7").Value = 1953
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1
$ws.Range("N7").Value = 2021
$ws.Range("O7").Value = 10
$ws.Range("P7").Value = 1

# J6 and J7 share a vertically-centered style with an explicit black font
# color. Build the combined format once in a scratch cell and copy it over
# (xlPasteFormats) so only a single new cell style is minted.
$scratch = $ws.Range("Z1")
$scratch.Font.Color = 0
$scratch.VerticalAlignment = -4108
[void]($scratch.Copy())
[void]($ws.Range("J6:J7").PasteSpecial(-4122))
[void]($scratch.Clear())
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Row 8 - BEIJING station (CHM00054511)
# ---------------------------------------------------------------------------
$ws.Range("E8").Value = "BEIJING"
$ws.Range("F8").Value = 39.93
$ws.Range("G8").Value = 116.28
$ws.Range("H8").Value = 55
$ws.Range("I8").Value = "CHM00054511"
$ws.Range("J8").Value = 54511
$ws.Range("K8").Value = 1951
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1
$ws.Range("N8").Value = 2021
$ws.Range("O8").Value = 9
$ws.Range("P8").Value = 29

# ---------------------------------------------------------------------------
# Row 9 - ANYANG station (CHM00053898)
# ---------------------------------------------------------------------------
$ws.Range("E9").Value = "ANYANG"
$ws.Range("F9").Value = 36.05
$ws.Range("G9").Value = 114.4
$ws.Range("H9").Value = 64
$ws.Range("I9").Value = "CHM00053898"
$ws.Range("J9").Value = 53898
$ws.Range("K9").Value = 1951
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1
$ws.Range("N9").Value = 2021
$ws.Range("O9").Value = 10
$ws.Range("P9").Value = 1

# ---------------------------------------------------------------------------
# Row 10 - HELLINIKON station (GR000016716)
# ---------------------------------------------------------------------------
$ws.Range("E10").Value = "HELLINIKON"
$ws.Range("F10").Value = 37.9
$ws.Range("G10").Value = 23.75
$ws.Range("H10").Value = 10
$ws.Range("I10").Value = "GR000016716"
$ws.Range("J10").Value = 16716
$ws.Range("K10").Value = 1955
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1
$ws.Range("N10").Value = 2021
$ws.Range("O10").Value = 10
$ws.Range("P10").Value = 1

# ---------------------------------------------------------------------------
# Row 11 - HELLINIKON station (GR000016717)
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = "HELLINIKON"
$ws.Range("F11").Value = 37.9
$ws.Range("G11").Value = 23.75
$ws.Range("H11").Value = 10
$ws.Range("I11").Value = "GR000016717"
$ws.Range("J11").Value = 16716
$ws.Range("K11").Value = 1955
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1
$ws.Range("N11").Value = 2021
$ws.Range("O11").Value = 10
$ws.Range("P11").Value = 1

# ---------------------------------------------------------------------------
# Row 12 - JUIZ_DE_FORA station (BR002143012); latitude sign corrected
# ---------------------------------------------------------------------------
$ws.Range("B12").Value = -21.23
$ws.Range("E12").Value = "JUIZ_DE_FORA"
$ws.Range("F12").Value = -21.77
$ws.Range("G12").Value = -43.35
$ws.Range("H12").Value = 911
$ws.Range("I12").Value = "BR002143012"
$ws.Range("J12").Value = 83692
$ws.Range("K12").Value = 1977
$ws.Range("L12").Value = 10
$ws.Range("M12").Value = 24
$ws.Range("N12").Value = 2017
$ws.Range("O12").Value = 8
$ws.Range("P12").Value = 7

# ---------------------------------------------------------------------------
# Row 13 - JUIZ_DE_FORA station (BR002143013); latitude sign corrected
# ---------------------------------------------------------------------------
$ws.Range("B13").Value = -21.23
$ws.Range("F13").Value = -21.77
$ws.Range("G13").Value = -43.35
$ws.Range("H13").Value = 911
$ws.Range("I13").Value = "BR002143013"
$ws.Range("J13").Value = 83692
$ws.Range("K13").Value = 1977
$ws.Range("L13").Value = 10
$ws.Range("M13").Value = 24
$ws.Range("N13").Value = 2017
$ws.Range("O13").Value = 8
$ws.Range("P13").Value = 7

# ---------------------------------------------------------------------------
# Row 14 - JUIZ_DE_FORA station (BR002143014); latitude sign corrected
# ---------------------------------------------------------------------------
$ws.Range("B14").Value = -21.23
$ws.Range("F14").Value = -21.77
$ws.Range("G14").Value = -43.35
$ws.Range("H14").Value = 911
$ws.Range("I14").Value = "BR002143014"
$ws.Range("J14").Value = 83692
$ws.Range("K14").Value = 1977
$ws.Range("L14").Value = 10
$ws.Range("M14").Value = 24
$ws.Range("N14").Value = 2017
$ws.Range("O14").Value = 8
$ws.Range("P14").Value = 7

# ---------------------------------------------------------------------------
# Fill the End_day (Q) formula down as a shared formula across Q3:Q14.
# ---------------------------------------------------------------------------
$ws.Range("Q3:Q14").Formula = "=N3-K3+1"

# ---------------------------------------------------------------------------
# Page setup: explicit portrait orientation.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# View tweaks: scroll the visible window and change the active selection.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 9
[void]($ws.Range("Q12:Q14").Select())
